$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The block of 8 rows that previously sat below the "ANA" (004210959) row
# (rows 8-15) is moved up to occupy rows 7-14, and the "ANA" row is moved
# to the end of that block (row 15) with its Saldo value changed to 16000.

$accounts = @("004526450","004432935","004498637","004693349","005064129","004473829","005684392","001000882","004210959")
$names    = @("MSD","JOSE","TIAGO","CATARINE","THIAGO","LUCAS","PAULO","AYRTON","ANA")
$saldos   = @(39000,32983.38,31324.89,30063.84,18817.87,18652.2,18065.7,16457.82,16000)

$startRow = 7
for ($i = 0; $i -lt $accounts.Length; $i++) {
    $r = $startRow + $i
    $acctCell = $ws.Cells.Item($r, 1)
    $acctCell.NumberFormat = "@"
    $acctCell.Value = $accounts[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $saldos[$i]
}
